# Update EMEP scaling mapping workbook to the new "year" sheet format:
#  - add 4 new columns (select_scaling_year, start_scaling_year, end_scaling_year, Comment)
#  - populate row 2 with real default values instead of placeholder "NA"s
#  - make the "year" sheet the active/selected tab (was "map")

$wb = $excel.ActiveWorkbook
$yearWs = $wb.Worksheets.Item("year")

# --- New header cells (row 1) ---
$yearWs.Range("E1").Value = "select_scaling_year"
$yearWs.Range("F1").Value = "start_scaling_year"
$yearWs.Range("G1").Value = "end_scaling_year"
$yearWs.Range("H1").Value = "Comment"

# --- New / updated data cells (row 2) ---
$yearWs.Range("A2").Value = "mkd"
$yearWs.Range("B2").Value = "all"
$yearWs.Range("E2").Value = "NA"
$yearWs.Range("F2").Value = 1990
$yearWs.Range("G2").Value = 2010
$yearWs.Range("H2").Value = "Don't scale to 1990-1991 drop so as to be closer to EMEP trend"

# --- Make "year" the selected/active sheet (moves tabSelected from "map") ---
$yearWs.Activate()
$yearWs.Rows("1:2").Select() | Out-Null
